# Fix: Elective lecture and tutorial scheduling
# Updates classroom/lab assignments in the Regular_Timetable, PreMid_Timetable
# and PostMid_Timetable sheets of the Semester 5 DSAI timetable workbook.

$wb = $excel.ActiveWorkbook

# Each of the three timetable sheets carries its own room assignment for the
# "DS302 / DS303 / CS307" block (rows 3,4,6,7,8,9) and the shared
# ELECTIVE_B4 / ELECTIVE_B5 basket table (rows 22-31). The room codes differ
# per sheet, the elective-basket edits are identical across all three sheets.

$sheetRoomMap = @{
    "Regular_Timetable" = @{ NewRoom = "C304"; NewLab = "L207" }
    "PreMid_Timetable"  = @{ NewRoom = "C305"; NewLab = "L207" }
    "PostMid_Timetable" = @{ NewRoom = "C102"; NewLab = "L207" }
}

foreach ($sheetName in $sheetRoomMap.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rooms = $sheetRoomMap[$sheetName]
    $newRoom = $rooms.NewRoom
    $newLab = $rooms.NewLab

    $ws.Range("B3").Value = "DS302 [$newRoom]"
    $ws.Range("C3").Value = "DS302 [$newRoom]"
    $ws.Range("D3").Value = "CS307 [$newRoom]"

    $ws.Range("B4").Value = "DS303 [$newRoom]"
    $ws.Range("C4").Value = "DS303 [$newRoom]"

    $ws.Range("B6").Value = "CS307 [$newRoom]"

    $ws.Range("B7").Value = "DS302 (Tutorial) [$newRoom]"

    $ws.Range("E8").Value = "CS307 (Lab) [$newLab]"

    $ws.Range("B9").Value = "DS303 (Tutorial) [$newRoom]"

    $ws.Range("E9").Value = "CS307 (Lab) [$newLab]"

    # ELECTIVE_B4 / ELECTIVE_B5 basket table: swap stale lab-building room
    # codes out of the lecture slot column and append the classroom code to
    # the previously room-less tutorial slot column.
    $ws.Range("D22").Value = "Tue 13:00-14:30 [C101], Thu 13:00-14:30 [C101]"
    $ws.Range("E22").Value = "Wed 14:30-15:30 [C101]"

    $ws.Range("D23").Value = "Tue 13:00-14:30 [C102], Thu 13:00-14:30 [C102]"
    $ws.Range("E23").Value = "Wed 14:30-15:30 [C102]"

    $ws.Range("D24").Value = "Tue 13:00-14:30 [C104], Thu 13:00-14:30 [C104]"
    $ws.Range("E24").Value = "Wed 14:30-15:30 [C104]"

    $ws.Range("D25").Value = "Tue 13:00-14:30 [C202], Thu 13:00-14:30 [C202]"
    $ws.Range("E25").Value = "Wed 14:30-15:30 [C202]"

    $ws.Range("D26").Value = "Mon 15:30-17:00 [C101], Wed 15:30-17:00 [C101]"
    $ws.Range("E26").Value = "Thu 14:30-15:30 [C101]"

    $ws.Range("D27").Value = "Mon 15:30-17:00 [C102], Wed 15:30-17:00 [C102]"
    $ws.Range("E27").Value = "Thu 14:30-15:30 [C102]"

    $ws.Range("D28").Value = "Mon 15:30-17:00 [C104], Wed 15:30-17:00 [C104]"
    $ws.Range("E28").Value = "Thu 14:30-15:30 [C104]"

    $ws.Range("D29").Value = "Mon 15:30-17:00 [C202], Wed 15:30-17:00 [C202]"
    $ws.Range("E29").Value = "Thu 14:30-15:30 [C202]"

    $ws.Range("D30").Value = "Mon 15:30-17:00 [C203], Wed 15:30-17:00 [C203]"
    $ws.Range("E30").Value = "Thu 14:30-15:30 [C203]"

    $ws.Range("D31").Value = "Mon 15:30-17:00 [C204], Wed 15:30-17:00 [C204]"
    $ws.Range("E31").Value = "Thu 14:30-15:30 [C204]"
}
